$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to be treated as text so that numeric-looking
# values (e.g. "1.000", "0.9996") are preserved as strings rather than converted
# to numbers, matching the original inline-string cell content.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Update Price (D) column values
$ws.Range("D2").Value = "29.479.78"
$ws.Range("D3").Value = "1.848.13"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D5").Value = "242.93"
$ws.Range("D6").Value = "0.6513"
$ws.Range("D7").Value = "1.000"
$ws.Range("D8").Value = "47.85"
$ws.Range("D9").Value = "0.07492"
$ws.Range("D11").Value = "24.46"
$ws.Range("D12").Value = "0.07626"
$ws.Range("D13").Value = "1.882.91"
$ws.Range("D14").Value = "5.014"
$ws.Range("D15").Value = "0.6839"
$ws.Range("D16").Value = "83.49"
$ws.Range("D17").Value = "0.000009446"
$ws.Range("D18").Value = "6.093"
$ws.Range("D19").Value = "29.520.96"
$ws.Range("D20").Value = "2.109.43"
$ws.Range("D21").Value = "236.96"
$ws.Range("D24").Value = "7.693"
$ws.Range("D26").Value = "157.41"
$ws.Range("D27").Value = "0.1417"
$ws.Range("D28").Value = "8.495"
$ws.Range("D29").Value = "17.80"
$ws.Range("D30").Value = "0.06066"
$ws.Range("D31").Value = "1.487"
$ws.Range("D33").Value = "4.134"
$ws.Range("D34").Value = "4.066"
$ws.Range("D35").Value = "1.181"
$ws.Range("D37").Value = "0.7232"
$ws.Range("D38").Value = "2.592"
$ws.Range("D39").Value = "2.796"
$ws.Range("D40").Value = "0.01780"
$ws.Range("D41").Value = "1.200.55"
$ws.Range("D42").Value = "6.230"
$ws.Range("D43").Value = "0.9068"
$ws.Range("D44").Value = "1.000"
$ws.Range("D45").Value = "2.022.09"
$ws.Range("D46").Value = "101.79"
$ws.Range("D47").Value = "66.50"
$ws.Range("D48").Value = "7.457"
$ws.Range("D51").Value = "9.139"

# Restore the default (Normal) style on the Price column now that the text values
# have been written, so no stray cell formatting is left behind.
$priceRange.Style = "Normal"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +3.80%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +4.52%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +11.13%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("E51").Value = "  -1.72%  "
